$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value2 = 738.2
$ws.Range("I43").Value2 = 366.33334
$ws.Range("J43").Value2 = 831.1667
$ws.Range("K43").Value2 = 366.33334
$ws.Range("L43").Value2 = 831.1667
$ws.Range("M43").Value2 = -297.33334
$ws.Range("N43").Value2 = -969.1667

# Hunk 1: ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value2 = 3508.3696
$ws.Range("I76").Value2 = 3326.5854
$ws.Range("K76").Value2 = 3326.5854
$ws.Range("M76").Value2 = -3011.5854

# Hunk 2: ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value2 = 3508.3696
$ws.Range("I79").Value2 = 3326.5854
$ws.Range("K79").Value2 = 3326.5854
$ws.Range("M79").Value2 = -2234.5854

# Hunk 3: ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value2 = 2985.182
$ws.Range("I106").Value2 = 3047.7144
$ws.Range("J106").Value2 = 2875.75
$ws.Range("K106").Value2 = 3047.7144
$ws.Range("L106").Value2 = 2875.75
$ws.Range("M106").Value2 = -2416.7144
$ws.Range("N106").Value2 = -4137.75

# Hunk 4: ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value2 = 1985.25
$ws.Range("I113").Value2 = 1622.45
$ws.Range("J113").Value2 = 2438.75
$ws.Range("K113").Value2 = 1622.45
$ws.Range("L113").Value2 = 2438.75
$ws.Range("M113").Value2 = 1631.55
$ws.Range("N113").Value2 = -8946.75

# Hunk 5: ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value2 = 2510.7026
$ws.Range("I137").Value2 = 2193.3225
$ws.Range("J137").Value2 = 4150.5
$ws.Range("K137").Value2 = 6579.967500000001
$ws.Range("L137").Value2 = 12451.5
$ws.Range("M137").Value2 = -4029.967500000001
$ws.Range("N137").Value2 = -17551.5

# Hunk 6: ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value2 = 5147.0713
$ws.Range("I138").Value2 = 1221.3103
$ws.Range("J138").Value2 = 7923.829
$ws.Range("K138").Value2 = 3663.9309
$ws.Range("L138").Value2 = 23771.487
$ws.Range("M138").Value2 = 1476.0691
$ws.Range("N138").Value2 = -34051.487

# Hunk 7: ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value2 = 1622.2084
$ws.Range("I132").Value2 = 1130.25
$ws.Range("J132").Value2 = 2606.125
$ws.Range("K132").Value2 = 3390.75
$ws.Range("L132").Value2 = 7818.375
$ws.Range("M132").Value2 = -860.75
$ws.Range("N132").Value2 = -12878.375

# Hunk 8: ARM row 138
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value2 = 38981.332
$ws.Range("J138").Value2 = 38981.332
$ws.Range("L138").Value2 = 38981.332
$ws.Range("N138").Value2 = -49261.332

# Hunk 9: BSM row 13
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value2 = 64450
$ws.Range("J13").Value2 = 64450
$ws.Range("L13").Value2 = 64450
$ws.Range("N13").Value2 = -64786

# Hunk 10: BSM row 50
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value2 = 0
$ws.Range("J50").Value2 = 0
$ws.Range("L50").Value2 = 0
$ws.Range("N50").ClearContents()

# Hunk 11: BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 3460.7
$ws.Range("I134").Value2 = 3484.9473
$ws.Range("K134").Value2 = 10454.8419
$ws.Range("M134").Value2 = -7919.841899999999

# Hunk 12: CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 158.57143
$ws.Range("I22").Value2 = 146
$ws.Range("K22").Value2 = 146
$ws.Range("M22").Value2 = 204

# Hunk 13: CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 5673.0938
$ws.Range("I31").Value2 = 8119.2666
$ws.Range("J31").Value2 = 3514.7058
$ws.Range("K31").Value2 = 8119.2666
$ws.Range("L31").Value2 = 3514.7058
$ws.Range("M31").Value2 = -7824.2666
$ws.Range("N31").Value2 = -4104.7058

# Hunk 14: CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value2 = 5673.0938
$ws.Range("I34").Value2 = 8119.2666
$ws.Range("J34").Value2 = 3514.7058
$ws.Range("K34").Value2 = 8119.2666
$ws.Range("L34").Value2 = 3514.7058
$ws.Range("M34").Value2 = -7917.2666
$ws.Range("N34").Value2 = -3918.7058

# Hunk 15: CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 2116478.2
$ws.Range("I58").Value2 = 3248207.8
$ws.Range("J58").Value2 = 3916.3333
$ws.Range("K58").Value2 = 3248207.8
$ws.Range("L58").Value2 = 3916.3333
$ws.Range("M58").Value2 = -3248004.8
$ws.Range("N58").Value2 = -4322.3333

# Hunk 16: CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value2 = 1888.36
$ws.Range("I99").Value2 = 1572.625
$ws.Range("J99").Value2 = 2449.6667
$ws.Range("K99").Value2 = 1572.625
$ws.Range("L99").Value2 = 2449.6667
$ws.Range("M99").Value2 = -74.625
$ws.Range("N99").Value2 = -5445.6667

# Hunk 17: CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value2 = 555.2941
$ws.Range("I107").Value2 = 432.83334
$ws.Range("J107").Value2 = 622.0909
$ws.Range("K107").Value2 = 432.83334
$ws.Range("L107").Value2 = 622.0909
$ws.Range("M107").Value2 = 1487.16666
$ws.Range("N107").Value2 = -4462.0909

# Hunk 18: CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value2 = 1888.36
$ws.Range("I126").Value2 = 1572.625
$ws.Range("J126").Value2 = 2449.6667
$ws.Range("K126").Value2 = 4717.875
$ws.Range("L126").Value2 = 7349.000100000001
$ws.Range("M126").Value2 = -2247.875
$ws.Range("N126").Value2 = -12289.0001

# Hunk 19: CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value2 = 2547.6
$ws.Range("I132").Value2 = 2137.0625
$ws.Range("J132").Value2 = 3277.4443
$ws.Range("K132").Value2 = 6411.1875
$ws.Range("L132").Value2 = 9832.332900000001
$ws.Range("M132").Value2 = -3881.1875
$ws.Range("N132").Value2 = -14892.3329

# Hunk 20: CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value2 = 3533.75
$ws.Range("I134").Value2 = 2534.8333
$ws.Range("J134").Value2 = 4033.2083
$ws.Range("K134").Value2 = 7604.499899999999
$ws.Range("L134").Value2 = 12099.6249
$ws.Range("M134").Value2 = -5069.499899999999
$ws.Range("N134").Value2 = -17169.6249

# Hunk 21: CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value2 = 2116478.2
$ws.Range("I136").Value2 = 3248207.8
$ws.Range("J136").Value2 = 3916.3333
$ws.Range("K136").Value2 = 9744623.399999999
$ws.Range("L136").Value2 = 11748.9999
$ws.Range("M136").Value2 = -9742073.399999999
$ws.Range("N136").Value2 = -16848.9999

# Hunk 22: CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value2 = 1800
$ws.Range("I17").Value2 = 800
$ws.Range("J17").Value2 = 2000
$ws.Range("K17").Value2 = 2400
$ws.Range("L17").Value2 = 6000
$ws.Range("M17").Value2 = -2231
$ws.Range("N17").Value2 = -6338

# Hunk 23: CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value2 = 9567
$ws.Range("J39").Value2 = 9567
$ws.Range("L39").Value2 = 28701
$ws.Range("N39").Value2 = -29289

# Hunk 24: CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value2 = 913.0278
$ws.Range("I122").Value2 = 390.83334
$ws.Range("J122").Value2 = 1017.4667
$ws.Range("K122").Value2 = 3517.50006
$ws.Range("L122").Value2 = 9157.2003
$ws.Range("M122").Value2 = -1067.50006
$ws.Range("N122").Value2 = -14057.2003

# Hunk 25: CUL row 126
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value2 = 3715.75
$ws.Range("I126").Value2 = 0
$ws.Range("J126").Value2 = 3715.75
$ws.Range("K126").Value2 = 0
$ws.Range("L126").Value2 = 11147.25
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value2 = -21027.25

# Hunk 26: GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value2 = 0
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = 0
$ws.Range("L5").Value2 = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

# Hunk 27: GSM row 9
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value2 = 6833.3335
$ws.Range("I9").Value2 = 800
$ws.Range("J9").Value2 = 9850
$ws.Range("K9").Value2 = 800
$ws.Range("L9").Value2 = 9850
$ws.Range("N9").Value2 = -10190
$ws.Range("M9").Value2 = -630

# Hunk 28: GSM row 53
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value2 = 19900
$ws.Range("I53").Value2 = 10000
$ws.Range("J53").Value2 = 29800
$ws.Range("K53").Value2 = 10000
$ws.Range("L53").Value2 = 29800
$ws.Range("M53").Value2 = -9369
$ws.Range("N53").Value2 = -31062

# Hunk 29: GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 3845.5908
$ws.Range("I122").Value2 = 4214.7334
$ws.Range("J122").Value2 = 3054.5715
$ws.Range("K122").Value2 = 12644.2002
$ws.Range("L122").Value2 = 9163.7145
$ws.Range("M122").Value2 = -10194.2002
$ws.Range("N122").Value2 = -14063.7145

# Hunk 30: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 2179.9688
$ws.Range("I132").Value2 = 2215.7368
$ws.Range("J132").Value2 = 2127.6924
$ws.Range("K132").Value2 = 6647.2104
$ws.Range("L132").Value2 = 6383.0772
$ws.Range("M132").Value2 = -4117.2104
$ws.Range("N132").Value2 = -11443.0772

# Hunk 31: LTW row 94
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value2 = 19460
$ws.Range("J94").Value2 = 19460
$ws.Range("L94").Value2 = 19460
$ws.Range("N94").Value2 = -20812

# Hunk 32: WVR row 104
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value2 = 0
$ws.Range("J104").Value2 = 0
$ws.Range("L104").Value2 = 0
$ws.Range("N104").ClearContents()

# Hunk 33: WVR row 118
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value2 = 67200
$ws.Range("J118").Value2 = 67200
$ws.Range("L118").Value2 = 67200
$ws.Range("N118").Value2 = -70514

# Hunk 34: WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 1825.6923
$ws.Range("I122").Value2 = 1346.4736
$ws.Range("K122").Value2 = 4039.4208
$ws.Range("M122").Value2 = -1589.4208

# Hunk 35: WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value2 = 1399.5714
$ws.Range("I126").Value2 = 1434.4706
$ws.Range("J126").Value2 = 1251.25
$ws.Range("K126").Value2 = 4303.4118
$ws.Range("L126").Value2 = 3753.75
$ws.Range("M126").Value2 = -1833.4118
$ws.Range("N126").Value2 = -8693.75

# Hunk 36: WVR row 130
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value2 = 23454.545
$ws.Range("J130").Value2 = 23454.545
$ws.Range("L130").Value2 = 23454.545
$ws.Range("N130").Value2 = -33494.545

# Hunk 37: WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 2312.3438
$ws.Range("I132").Value2 = 1060.8
$ws.Range("J132").Value2 = 3416.647
$ws.Range("K132").Value2 = 3182.4
$ws.Range("L132").Value2 = 10249.941
$ws.Range("M132").Value2 = -652.3999999999996
$ws.Range("N132").Value2 = -15309.941
